# Updated cryptos list on Fri Jan 26 11:40:01 UTC 2024 with GitHub Actions
# Refreshes the per-coin Price (D) / Volume(1h) (E) columns with the latest
# scrape, plus the Celestia/Stellar and VeChain/FraxShare row swaps (name,
# link, price and volume all move together) captured by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed prices are plain decimals (e.g. "1.00", "91.80") that Excel's
# COM layer would otherwise auto-convert to numbers (dropping the trailing
# zero / exact text). Force those Price cells to Text first so the literal
# string from the source feed is preserved, matching the existing sheet
# (every Price/Volume(1h) cell is stored as text).
$numericLookingPriceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D16",
    "D18",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49"
)

foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.187.23'
$ws.Range("E2").Value = '  +2.26%  '
$ws.Range("D3").Value = '2.255.47'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '302.94'
$ws.Range("E5").Value = '  +2.90%  '
$ws.Range("D6").Value = '91.80'
$ws.Range("E6").Value = '  +3.17%  '
$ws.Range("D7").Value = '0.521'
$ws.Range("E7").Value = '  +1.40%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '32.22'
$ws.Range("E10").Value = '  +5.45%  '
$ws.Range("D11").Value = '52.75'
$ws.Range("E11").Value = '  +6.08%  '
$ws.Range("D12").Value = '0.0794'
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").Value = '  +3.09%  '
$ws.Range("D14").Value = '6.59'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '2.597.39'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '14.16'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '2.255.95'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '0.753'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = '41.099.20'
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("D20").Value = '11.89'
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").Value = '0.0₃0903'
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").Value = '5.86'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '66.69'
$ws.Range("E23").Value = '  +1.40%  '
$ws.Range("D24").Value = '240.47'
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  +3.72%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '1.88'
$ws.Range("E27").Value = '  +2.58%  '
$ws.Range("D28").Value = '24.05'
$ws.Range("E28").Value = '  +5.08%  '
$ws.Range("E29").Value = '  -2.70%  '
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("D31").Value = '158.17'
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("D32").Value = '33.50'
$ws.Range("E32").Value = '  +3.47%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '5.11'
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("D35").Value = '3.07'
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").Value = '0.0730'
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("E37").Value = '  +6.90%  '
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = '0.115'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '16.40'
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("E41").Value = '  +5.26%  '
$ws.Range("D42").Value = '3.92'
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").Value = '2.091.50'
$ws.Range("E43").Value = '  -2.12%  '
$ws.Range("D44").Value = '20.10'
$ws.Range("E44").Value = '  +9.99%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '10.31'
$ws.Range("E45").Value = '  +4.94%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0277'
$ws.Range("E46").Value = '  +2.59%  '
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").Value = '  +8.00%  '
$ws.Range("D48").Value = '1.86'
$ws.Range("E48").Value = '  -13.38%  '
$ws.Range("D49").Value = '1.53'
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("D50").Value = '2.469.09'
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("E51").Value = '  +3.48%  '
